$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 416000000.0
$ws.Range("C4").Value = 405000000.0
$ws.Range("D4").Value = 387000000.0
$ws.Range("E4").Value = 382000000.0
$ws.Range("F4").Value = 363000000.0

$ws.Range("B12").Value = 1602000000.0
$ws.Range("C12").Value = 1980000000.0
$ws.Range("D12").Value = 1691000000.0
$ws.Range("E12").Value = 1657000000.0
$ws.Range("F12").Value = 1454000000.0

$ws.Range("B21").Value = 5474000000.0
$ws.Range("C21").Value = 5368000000.0
$ws.Range("D21").Value = 5161000000.0
$ws.Range("F21").Value = 5173000000.0

$ws.Range("G35").Value = 18825000000.0
$ws.Range("G36").Value = 18893000000.0
